# Latest Eric Edits Made
# Change the two oval (ellipse) shapes into rectangles, renaming them to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$msoShapeRectangle = 1

$oval2 = $s.Shapes.Item(2)
$oval2.AutoShapeType = $msoShapeRectangle
$oval2.Name = "PPRect#2"

$oval3 = $s.Shapes.Item(3)
$oval3.AutoShapeType = $msoShapeRectangle
$oval3.Name = "PPRect#3"
